$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6822801232337952
$ws.Range("B1").Value = 3.386894464492798
$ws.Range("C1").Value = 3.675288915634155
$ws.Range("D1").Value = 1.097281932830811
$ws.Range("E1").Value = 0.9883050322532654
